# Update "top-level-description" (column B) cells with revised, more
# concise wording for several rows on the Specification sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Specification")

$updates = @{
    21  = "Name and contact information if an agent is being used."
    25  = "Name and contact information if an agent is being used."
    33  = "Telephone number and email address of the applicant."
    37  = "Name and contact information for the parties making the application."
    43  = "Checking whether all the requirements of the form have been met, such as proof of payment or supporting documentation."
    44  = "What community consultation activities have taken place as part of the application"
    46  = "Details of any conflict of interest that may exist between the applicant and planning authority."
    49  = "Signed and dated verification of the application's accuracy."
    52  = "Details of any demolition that needs to take place as part of the development proposal."
    61  = "Whether the applicant has obtained a Certificate of Immunity (COI) meaning the building in question cannot be listed"
    63  = "Details of any changes being made to a listed building as part of development works"
    67  = "The grade of any listed building affected by the proposed development."
    70  = "What materials are being used for the proposed development"
    78  = "Who will be affected by the proposal and whether they have been notified, such as agricultural tenants"
    92  = "Details of pre-application advice received from the planning authority"
    97  = "What development, works or change of use is proposed"
    102 = "Details of any other development proposals made for the site"
    106 = "Where the proposed development will be built."
    115 = "Information to help the planning authority arrange a site visit"
}

foreach ($row in $updates.Keys) {
    $ws.Range("B$row").Value = $updates[$row]
}
